$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price (D) and 1h volume change (E) figures
# D-column cells are forced to Text format before assignment so that
# values such as "102.60" or "1.20" are not silently re-interpreted
# as numbers (which would drop the significant trailing zero).

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '38.727.60'
$ws.Range('E2').Value = '  +0.17%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.098.19'
$ws.Range('E3').Value = '  +0.23%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.59'
$ws.Range('E5').Value = '  -0.57%  '

$ws.Range('E6').Value = '  +0.37%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '62.07'
$ws.Range('E7').Value = '  +1.40%  '

$ws.Range('E9').Value = '  +1.87%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0839'
$ws.Range('E10').Value = '  -0.42%  '

$ws.Range('E11').Value = '  -1.22%  '

$ws.Range('E12').Value = '  +5.43%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '2.409.80'
$ws.Range('E13').Value = '  +0.52%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '22.03'
$ws.Range('E14').Value = '  -1.57%  '

$ws.Range('E15').Value = '  +3.48%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.53'
$ws.Range('E16').Value = '  +1.15%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.140.37'
$ws.Range('E17').Value = '  +2.82%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '38.700.35'
$ws.Range('E18').Value = '  +0.36%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '71.58'
$ws.Range('E19').Value = '  +0.79%  '

$ws.Range('E20').Value = '  +0.33%  '

$ws.Range('E21').Value = '  +0.37%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '227.83'
$ws.Range('E22').Value = '  +0.66%  '

$ws.Range('E23').Value = '  -0.01%  '

$ws.Range('E24').Value = '  -2.01%  '

$ws.Range('E25').Value = '  -0.60%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.62'
$ws.Range('E26').Value = '  +1.78%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '171.49'
$ws.Range('E27').Value = '  +0.64%  '

$ws.Range('E28').Value = '  +2.13%  '

$ws.Range('E29').Value = '  +3.62%  '

$ws.Range('E30').Value = '  +0.75%  '

$ws.Range('E31').Value = '  +8.15%  '

$ws.Range('E32').Value = '  +0.06%  '

$ws.Range('E33').Value = '  +0.84%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.73'
$ws.Range('E34').Value = '  -0.75%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '7.02'
$ws.Range('E35').Value = '  +7.14%  '

$ws.Range('E36').Value = '  +1.63%  '

$ws.Range('E37').Value = '  +0.04%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.54'
$ws.Range('E38').Value = '  -0.82%  '

$ws.Range('E39').Value = '  +0.07%  '

$ws.Range('E40').Value = '  -2.36%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '102.60'
$ws.Range('E41').Value = '  +2.41%  '

$ws.Range('E42').Value = '  +2.84%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.526.79'
$ws.Range('E43').Value = '  -1.26%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.20'
$ws.Range('E44').Value = '  +6.64%  '

$ws.Range('E45').Value = '  -0.92%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.81'
$ws.Range('E46').Value = '  +2.11%  '

$ws.Range('E47').Value = '  -0.79%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '4.15'
$ws.Range('E48').Value = '  -0.34%  '

$ws.Range('E49').Value = '  +1.61%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.296.10'
$ws.Range('E51').Value = '  +0.31%  '

Write-Output "Updated 71 cells across 48 rows"
